$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeouts) values replacing the old Strike# values in column G,
# keyed by row number (row 2 = first data row).
$kValues = @{
    2  = 3
    3  = 2
    4  = 2
    6  = 2
    7  = 0
    8  = 1
    9  = 4
    10 = 2
    11 = 2
    12 = 5
    13 = 0
    14 = 1
    15 = 2
    16 = 0
    17 = 0
    18 = 1
    19 = 4
    20 = 1
    21 = 0
    22 = 2
    23 = 3
    24 = 2
    25 = 2
    26 = 1
    27 = 0
    28 = 3
    29 = 2
    30 = 2
    31 = 1
    32 = 4
    33 = 2
    34 = 2
    35 = 3
    36 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
